$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose text would otherwise be auto-coerced to a number by Excel's
# type inference are written with an explicit Text number format, then
# restored to the default "Normal" style so the cell keeps looking like
# every other untouched cell once the text has been committed.

$ws.Range("D2").Value = "22.451.43"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").Value = "1.567.96"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "288.29"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.62%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3725"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.78%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "48.25"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3321"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07472"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.129"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.86%  "
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.76"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.962"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.901"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.33%  "
$ws.Range("D16").Value = "1.568.31"
$ws.Range("E16").Value = "  -0.15%  "
$ws.Range("E17").Value = "  -0.83%  "
$ws.Range("E18").Value = "  -0.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "87.92"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.90%  "
$ws.Range("E20").Value = "  +0.15%  "
$ws.Range("E21").Value = "  -0.28%  "
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.98%  "
$ws.Range("D24").Value = "22.449.88"
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("E25").Value = "  +1.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.560"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.44%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.45"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.017"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "123.99"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.86%  "
$ws.Range("D31").Value = "1.743.65"
$ws.Range("E31").Value = "  -0.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.055"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.87%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.128"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.78%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.619"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08277"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02453"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2274"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.27%  "
$ws.Range("E39").Value = "  -3.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.353"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.288"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.60%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6277"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.24"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.95%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.002"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.83"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.71%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6127"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.778"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.043"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "125.42"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.97%  "
$ws.Range("E50").Value = "  -2.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07239"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.85%  "
